$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VIDEOS")

# New "Rounded" column (K) header
$ws.Range("K1").Value = "Rounded"

# Per-video rounded-duration values (K5:K204), one per row matching A5:A204 (1..200)
$roundedValues = @(35,49,49,49,52,51,65,59,71,55,83,54,5,42,49,58,57,46,59,3,57,20,24,18,24,20,13,8,28,26,19,15,23,18,16,25,20,27,19,23,17,16,30,29,21,23,32,23,28,16,24,33,19,16,25,18,24,22,24,20,27,14,21,30,18,27,28,29,24,36,30,6,33,26,23,24,17,21,19,20,22,22,27,23,13,19,18,26,25,25,30,27,34,24,19,23,21,22,18,27,26,17,19,16,21,20,15,24,15,23,15,22,19,16,18,18,24,12,20,11,19,14,20,19,20,9,25,20,19,11,17,17,23,25,10,26,24,31,10,25,14,18,29,18,14,17,20,21,13,23,19,19,21,20,19,21,11,20,22,8,29,20,14,11,29,30,18,15,8,23,25,8,18,23,21,10,8,27,22,22,22,16,21,13,20,12,19,10,5,25,21,19,15,16,22,20,11,17,14,24)

for ($i = 0; $i -lt $roundedValues.Length; $i++) {
    $row = 5 + $i
    $ws.Cells.Item($row, 11).Value = $roundedValues[$i]
}

# Status corrections: rows 13 & 14 moved to Complete, row 15 moved to In Progress
$ws.Range("H13").Value = "Complete"
$ws.Range("H14").Value = "Complete"
$ws.Range("H15").Value = "In Progress"

# Summary block under the table (rows 206-208)
$ws.Range("D206").Value = 280734
$ws.Range("D207").Value = 4679
$ws.Range("D208").Value = "78 hr"

$ws.Range("G206").Value = "Progress:"
$ws.Range("H206").Formula = '=SUMIFS(K5:K204, H5:H204, "Complete")'
$ws.Range("K206").Formula = "=SUM(K5:K204)"
